$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new column before the existing "area" column (C), shifting the
# old C:G (area/1500, owner, register_date, register_reason, acquire_value)
# one column to the right to D:H. Insert() copies the left neighbour's
# formatting onto the new column automatically.
$ws.Columns.Item(3).Insert()

# Extend the header row (styled like B1, bold+border) and the data row
# (styled like B2, plain) out to column N so the new cells pick up
# matching formatting.
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)

# --- Row 1 (header labels) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 (data values) ---
$ws.Range("B2").Value = "HONDACITY"
$ws.Range("C2").Value = 1500
$ws.Range("D2").Value = "王育敏"
$ws.Range("E2").Value = "87年10月01日"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = "400000(超過五年）"
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# Leading apostrophe forces this to be stored as text rather than being
# auto-parsed into a date serial number.
$ws.Range("J2").Value = "'2012-04-12"
$ws.Range("K2").Value = "王育敏"
$ws.Range("L2").Value = 1728
$ws.Range("M2").Value = "tmp48bc1"
$ws.Range("N2").Value = 29

"done"
